# Update "南宁-漫展信息" workbook: add a new exhibition event
#   "南宁·0316全职only-全明星周末" (2024-03-16) into the 展览 (sheet1) and
#   全部类型 (sheet4) sheets, right before the existing "南宁·草莓动漫节"
#   row, and bump a couple of "想去人数" (interest count) numbers that
#   were refreshed at the same time (F2 on every touched sheet, plus the
#   rows that got pushed down one slot).

$wb = $excel.ActiveWorkbook

function Update-Sheet($ws, $hasKanongRow) {

    # The overall "想去人数" figure for the first (冬典) event ticked up.
    $ws.Range("F2").Value = 9264

    if ($hasKanongRow) {
        # 全部类型 sheet currently has 6 data rows (rows 2-6); shift the
        # bottom two down by one to make room for the new row at 4.
        $ws.Range("A6:I6").Copy($ws.Range("A7:I7"))
        $ws.Range("A5:I5").Copy($ws.Range("A6:I6"))
        $ws.Range("A4:I4").Copy($ws.Range("A5:I5"))

        # Row 6 is the 卡农 row that just moved down from row 5 -> 6;
        # its running index is unchanged (still the 5th entry).
        $ws.Range("A6").Value = 5

        # Row 7 is the DACG row that just moved down from row 6 -> 7;
        # its running index and "想去人数" also refreshed.
        $ws.Range("A7").Value = 6
        $ws.Range("F7").Value = 458
    } else {
        # 展览 sheet currently has 5 data rows (rows 2-5); shift the
        # bottom row down by one to make room for the new row at 4.
        $ws.Range("A5:I5").Copy($ws.Range("A6:I6"))
        $ws.Range("A4:I4").Copy($ws.Range("A5:I5"))

        # Row 6 is the DACG row that just moved down from row 5 -> 6;
        # its "想去人数" also refreshed.
        $ws.Range("A6").Value = 5
        $ws.Range("F6").Value = 458
    }

    # Row 5 is now the 草莓动漫节 row that moved down from row 4 -> 5;
    # its running index and "想去人数" also refreshed.
    $ws.Range("A5").Value = 4
    $ws.Range("F5").Value = 495

    # Row 4 becomes the brand new event, overwriting the old
    # 草莓动漫节 data that was copied down to row 5 already.
    $ws.Range("C4").Value = "南宁·0316全职only-全明星周末"
    $ws.Range("D4").Value = "大学东路158号 维也纳国际酒店"
    $ws.Range("E4").Value = "2024.03.16 10:30-03.16 17:00"
    $ws.Range("F4").Value = 3
    $ws.Range("G4").Value = 65
    $ws.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=81834"
    $ws.Range("I4").Value = "//i2.hdslb.com/bfs/openplatform/202402/5RPyTNNO1707363370492.jpeg"
}

foreach ($ws in $wb.Worksheets) {
    if ($ws.Name -eq "展览") {
        Update-Sheet $ws $false
    } elseif ($ws.Name -eq "全部类型") {
        Update-Sheet $ws $true
    }
}
